$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.104.42"
$ws.Range("E2").Value = "  +0.12%  "

$ws.Range("D3").Value = "3.561.57"
$ws.Range("E3").Value = "  +2.39%  "

$ws.Range("E4").Value = "  -0.03%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "606.19"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.49%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "145.16"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.59%  "

$ws.Range("D7").Value = "3.559.34"
$ws.Range("E7").Value = "  +2.33%  "

$ws.Range("E8").Value = "  +0.43%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.491"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +3.49%  "

$ws.Range("E10").Value = "  +1.26%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "7.93"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -3.39%  "

$ws.Range("E12").Value = "  +0.13%  "

$ws.Range("D13").Value = "4.162.47"
$ws.Range("E13").Value = "  +2.34%  "

$ws.Range("E14").Value = "  +2.03%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "29.99"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.12%  "

$ws.Range("D16").Value = "3.558.77"
$ws.Range("E16").Value = "  +2.28%  "

$ws.Range("D17").Value = "66.201.28"
$ws.Range("E17").Value = "  +0.13%  "

$ws.Range("E18").Value = "  -0.98%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.42"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +9.74%  "

$ws.Range("E20").Value = "  +0.86%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.85"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.96%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "429.51"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +2.26%  "

$ws.Range("E23").Value = "  +4.70%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "79.19"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +2.22%  "

$ws.Range("D25").Value = "3.700.92"
$ws.Range("E25").Value = "  +2.24%  "

$ws.Range("E26").Value = "  +0.11%  "

$ws.Range("E27").Value = "  +2.15%  "

$ws.Range("E28").Value = "  +1.84%  "

$ws.Range("E29").Value = "  +0.10%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "9.11"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -3.05%  "

$ws.Range("E31").Value = "  -0.04%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "25.56"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.94%  "

$ws.Range("E33").Value = "  -1.31%  "

$ws.Range("D34").Value = "3.555.57"
$ws.Range("E34").Value = "  +2.35%  "

$ws.Range("E35").Value = "  -6.24%  "

$ws.Range("E37").Value = "  +1.96%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "7.86"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +3.19%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.61"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.15%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "174.74"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +2.92%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0848"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.12%  "

$ws.Range("E43").Value = "  +2.04%  "

$ws.Range("E44").Value = "  +0.83%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.94"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.94%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "46.13"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.06%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "25.83"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.14%  "

$ws.Range("E48").Value = "  +0.06%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.38"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.55%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "23.50"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +9.80%  "

$ws.Range("E51").Value = "  +0.05%  "
